# Update "想去人数" (want-to-go count) figures in the F column on both the
# "展览" sheet and the "全部类型" sheet, which mirrors the same data.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 3001
    "F5"  = 6756
    "F6"  = 1746
    "F7"  = 23
    "F9"  = 61
    "F10" = 123
    "F11" = 10
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
